# Insert a new first column ("ID") before the existing data, shifting the
# current columns A:E (now containing headers A, B, C, D, F) to B:F, and
# populate the new column with the sample identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns one place to the right.
$ws.Columns.Item(1).Insert()

# Copy the header formatting (bold font, border, centered alignment) from
# the neighboring header cell onto the new header cell.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Header text for the newly inserted column.
$ws.Range("A1").Value = "ID"

# Row identifiers for each data row (rows 2-25).
$ids = @(
    "Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95",
    "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22",
    "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16"
)

for ($i = 0; $i -lt $ids.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
